$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated F-column (movement rate) values for most rows ---
# (row 15 keeps its original formula =187.8826/2 and is left untouched)
$ws.Range("F2").Value  = 62.692799999999998
$ws.Range("F3").Value  = 84.306510000000003
$ws.Range("F4").Value  = 110.1503
$ws.Range("F5").Value  = 79.866299999999995
$ws.Range("F6").Value  = 79.421049999999994
$ws.Range("F7").Value  = 83.33869
$ws.Range("F8").Value  = 87.521619999999999
$ws.Range("F9").Value  = 92.282579999999996
$ws.Range("F10").Value = 119.7243
$ws.Range("F11").Value = 95.382819999999995
$ws.Range("F12").Value = 88.216070000000002
$ws.Range("F13").Value = 64.697630000000004
$ws.Range("F14").Value = 96.015469999999993
# F16 previously held the erroneous formula =113.0523/2 -- replace with a plain value
$ws.Range("F16").Value = 85.791179999999997
$ws.Range("F17").Value = 73.799009999999996
$ws.Range("F18").Value = 85.265659999999997
$ws.Range("F19").Value = 58.121429999999997

# --- G column: daily distance = F * 24 ---
$ws.Range("G2").Formula = "=F2*24"
$ws.Range("G3:G19").Formula = "=F3*24"

# --- Updated L-column values for most rows ---
# (row 6 keeps its original value unchanged)
$ws.Range("L2").Value  = 58.544849999999997
$ws.Range("L3").Value  = 50.364069999999998
$ws.Range("L4").Value  = 59.245980000000003
$ws.Range("L5").Value  = 64.57208
$ws.Range("L7").Value  = 67.105860000000007
$ws.Range("L8").Value  = 57.643569999999997
$ws.Range("L9").Value  = 53.38297
$ws.Range("L10").Value = 68.015469999999993
$ws.Range("L11").Value = 181.22890000000001
$ws.Range("L12").Value = 52.305050000000001
$ws.Range("L13").Value = 43.635570000000001
$ws.Range("L14").Value = 66.867980000000003
$ws.Range("L15").Value = 85.791179999999997
$ws.Range("L16").Value = 79.165719999999993
$ws.Range("L17").Value = 73.706360000000004
$ws.Range("L18").Value = 79.614689999999996
$ws.Range("L19").Value = 46.320880000000002

# --- M column: daily distance = L * 24 ---
$ws.Range("M2").Formula = "=L2*24"
$ws.Range("M3:M19").Formula = "=L3*24"

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 8.5
$ws.Columns.Item(8).ColumnWidth = 13.33

# --- Selection state ---
$ws.Range("I16").Select() | Out-Null
